$d = $word.ActiveDocument

# --- 1. Date: "May 18, 2022" -> "May 19, 2022" (the posting date moves a day) ---
$p1 = $d.Paragraphs(3).Range
$found1 = $p1.Find.Execute('May 18, 2022', $true, $false, $false, $false, $false, $true, 1, $false, 'May 19, 2022', 2)
Write-Host "1) date updated:" $found1

# --- 2. Append " (fork it first if you want)" right after the repo link ---
$p2 = $d.Paragraphs(57).Range
$p2.MoveEnd(1, -1) | Out-Null          # exclude the paragraph mark
$p2.Collapse(0) | Out-Null             # collapse to the very end of the paragraph text
$p2.InsertAfter(' (fork it first if you want)')
Write-Host "2) after:" $d.Paragraphs(57).Range.Text

# --- 3. "fork the repo and name the new fork after yourself somehow"
#        -> "checkout a new branch and rename it after yourself somehow" ---
$p3 = $d.Paragraphs(58).Range
$found3 = $p3.Find.Execute('fork the repo and name the new fork after yourself somehow', $true, $false, $false, $false, $false, $true, 1, $false, 'checkout a new branch and rename it after yourself somehow', 2)
Write-Host "3) branch instructions updated:" $found3

# --- 4. "...file in your fork" -> "...file in your branch" ---
$p4 = $d.Paragraphs(59).Range
$found4 = $p4.Find.Execute('/` file in your fork', $true, $false, $false, $false, $false, $true, 1, $false, '/` file in your branch', 2)
Write-Host "4) upload instructions updated:" $found4

# --- 5. Prepend "Git ships with RStudio and " before "Happy with Git..." ---
$p5 = $d.Paragraphs(61).Range
$found5 = $p5.Find.Execute('Happy with Git is an exception', $true, $false, $false, $false, $false, $true, 1, $false, 'Git ships with RStudio and Happy with Git is an exception', 2)
Write-Host "5) RStudio intro added:" $found5

# --- 6. " resource for getting going with Git.  " -> "...going with Rstudio and Git.  " ---
$p6 = $d.Paragraphs(61).Range
$found6 = $p6.Find.Execute('resource for getting going with Git.', $true, $false, $false, $false, $false, $true, 1, $false, 'resource for getting going with Rstudio and Git.', 2)
Write-Host "6) Rstudio mention added:" $found6

Write-Host "Final paragraph 61:" $d.Paragraphs(61).Range.Text
